$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook lists IATI sector codes together with their parent
# "category" and "group" metadata in columns D (category-name),
# E (category-code), F (group-name) and G (group-code).
#
# The upstream data generator changed the column order so that the
# group-code now comes first, i.e. for every row (including the header)
# the four values in D:G are rotated one position to the right:
#   new D = old G
#   new E = old D
#   new F = old E
#   new G = old F

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    $d = $dCell.Value2
    $e = $eCell.Value2
    $f = $fCell.Value2
    $g = $gCell.Value2

    $dCell.Value = $g
    $eCell.Value = $d
    $fCell.Value = $e
    $gCell.Value = $f
}
